$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Update the SqlIP value in E2 to a new, distinct IP address
$ws.Range("E2").Value = "192.168.0.24"

# Update the active selection to H6 (matches the diff's sheetView selection)
$ws.Range("H6").Select()
